$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 775.8333
$ws.Range("I19").Value = 1035
$ws.Range("J19").Value = 516.6667
$ws.Range("K19").Value = 1035
$ws.Range("L19").Value = 516.6667
$ws.Range("M19").Value = -860
$ws.Range("N19").Value = -866.6667
# Row 57
$ws.Range("H57").Value = 23780
$ws.Range("J57").Value = 23780
$ws.Range("L57").Value = 71340
$ws.Range("N57").Value = -72338
# Row 105
$ws.Range("H105").Value = 49667
$ws.Range("J105").Value = 49667
$ws.Range("L105").Value = 49667
$ws.Range("N105").Value = -56655
# Row 137
$ws.Range("H137").Value = 3055.0247
$ws.Range("I137").Value = 1284.5714
$ws.Range("K137").Value = 3853.7142
$ws.Range("M137").Value = -1303.7142

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1852.3704
$ws.Range("I2").Value = 1860.3478
$ws.Range("J2").Value = 1806.5
$ws.Range("K2").Value = 1860.3478
$ws.Range("L2").Value = 1806.5
$ws.Range("M2").Value = -1747.3478
$ws.Range("N2").Value = -2032.5
# Row 32
$ws.Range("H32").Value = 22352.436
$ws.Range("I32").Value = 21681.434
$ws.Range("J32").Value = 37450
$ws.Range("K32").Value = 21681.434
$ws.Range("L32").Value = 37450
$ws.Range("M32").Value = -21394.434
$ws.Range("N32").Value = -38024
# Row 45
$ws.Range("H45").Value = 58825730
$ws.Range("I45").Value = 71430550
$ws.Range("J45").Value = 3238
$ws.Range("K45").Value = 71430550
$ws.Range("L45").Value = 3238
$ws.Range("M45").Value = -71430173
$ws.Range("N45").Value = -3992
# Row 110
$ws.Range("H110").Value = 1937.7307
$ws.Range("I110").Value = 1909.1177
$ws.Range("J110").Value = 1991.7778
$ws.Range("K110").Value = 1909.1177
$ws.Range("L110").Value = 1991.7778
$ws.Range("M110").Value = 135.8823
$ws.Range("N110").Value = -6081.7778
# Row 116
$ws.Range("H116").Value = 1852.3704
$ws.Range("I116").Value = 1860.3478
$ws.Range("J116").Value = 1806.5
$ws.Range("K116").Value = 1860.3478
$ws.Range("L116").Value = 1806.5
$ws.Range("M116").Value = 433.6522
$ws.Range("N116").Value = -6394.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1852.3704
$ws.Range("I3").Value = 1860.3478
$ws.Range("J3").Value = 1806.5
$ws.Range("K3").Value = 1860.3478
$ws.Range("L3").Value = 1806.5
$ws.Range("M3").Value = -1746.3478
$ws.Range("N3").Value = -2034.5
# Row 96
$ws.Range("H96").Value = 13076.4
$ws.Range("I96").Value = 5382
$ws.Range("K96").Value = 5382
$ws.Range("M96").Value = -2636
# Row 130
$ws.Range("H130").Value = 53493
$ws.Range("J130").Value = 53493
$ws.Range("L130").Value = 53493
$ws.Range("N130").Value = -63533

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2506
$ws.Range("I16").Value = 2505.5
$ws.Range("J16").Value = 2506.5
$ws.Range("K16").Value = 2505.5
$ws.Range("L16").Value = 2506.5
$ws.Range("M16").Value = -2218.5
$ws.Range("N16").Value = -3080.5
# Row 31
$ws.Range("H31").Value = 5500.523
$ws.Range("I31").Value = 1923.3334
$ws.Range("J31").Value = 9793.15
$ws.Range("K31").Value = 1923.3334
$ws.Range("L31").Value = 9793.15
$ws.Range("M31").Value = -1628.3334
$ws.Range("N31").Value = -10383.15
# Row 32
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -684
$ws.Range("N32").ClearContents()
# Row 34
$ws.Range("H34").Value = 5500.523
$ws.Range("I34").Value = 1923.3334
$ws.Range("J34").Value = 9793.15
$ws.Range("K34").Value = 1923.3334
$ws.Range("L34").Value = 9793.15
$ws.Range("M34").Value = -1721.3334
$ws.Range("N34").Value = -10197.15
# Row 52
$ws.Range("H52").Value = 40911.332
$ws.Range("J52").Value = 40911.332
$ws.Range("L52").Value = 40911.332
$ws.Range("N52").Value = -41499.332
# Row 113
$ws.Range("H113").Value = 2506
$ws.Range("I113").Value = 2505.5
$ws.Range("J113").Value = 2506.5
$ws.Range("K113").Value = 2505.5
$ws.Range("L113").Value = 2506.5
$ws.Range("M113").Value = -335.5
$ws.Range("N113").Value = -6846.5
# Row 139
$ws.Range("H139").Value = 49563.8
$ws.Range("J139").Value = 50404.223
$ws.Range("L139").Value = 50404.223
$ws.Range("N139").Value = -60684.223
# Row 141
$ws.Range("H141").Value = 25498.5
$ws.Range("I141").Value = 18500
$ws.Range("J141").Value = 28997.75
$ws.Range("K141").Value = 18500
$ws.Range("L141").Value = 28997.75
$ws.Range("M141").Value = -13320
$ws.Range("N141").Value = -39357.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 166666960
$ws.Range("J80").Value = 166666960
$ws.Range("L80").Value = 500000880
$ws.Range("N80").Value = -500002752
# Row 83
$ws.Range("H83").Value = 166666960
$ws.Range("J83").Value = 166666960
$ws.Range("L83").Value = 1500002640
$ws.Range("N83").Value = -1500012000
# Row 92
$ws.Range("H92").Value = 789.7273
$ws.Range("I92").Value = 672
$ws.Range("J92").Value = 995.75
$ws.Range("K92").Value = 2016
$ws.Range("L92").Value = 2987.25
$ws.Range("M92").Value = -768
$ws.Range("N92").Value = -5483.25
# Row 95
$ws.Range("H95").Value = 2700000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 2700000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 8100000
$ws.Range("N95").Value = -8104118
$ws.Range("M95").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 101
$ws.Range("H101").Value = 48246
$ws.Range("J101").Value = 48246
$ws.Range("L101").Value = 48246
$ws.Range("N101").Value = -54736
# Row 108
$ws.Range("H108").Value = 20114
$ws.Range("J108").Value = 20114
$ws.Range("L108").Value = 20114
$ws.Range("N108").Value = -27794
# Row 122
$ws.Range("H122").Value = 1439.25
$ws.Range("I122").Value = 1465.8
$ws.Range("J122").Value = 1359.6
$ws.Range("K122").Value = 4397.4
$ws.Range("L122").Value = 4078.8
$ws.Range("M122").Value = -1947.4
$ws.Range("N122").Value = -8978.799999999999
# Row 126
$ws.Range("H126").Value = 10584.521
$ws.Range("I126").Value = 17080.309
$ws.Range("J126").Value = 2140
$ws.Range("K126").Value = 51240.927
$ws.Range("L126").Value = 6420
$ws.Range("M126").Value = -48770.927
$ws.Range("N126").Value = -11360
# Row 129
$ws.Range("H129").Value = 31499.8
$ws.Range("J129").Value = 31499.8
$ws.Range("L129").Value = 31499.8
$ws.Range("N129").Value = -41499.8
# Row 137
$ws.Range("H137").Value = 50282
$ws.Range("J137").Value = 50282
$ws.Range("L137").Value = 50282
$ws.Range("N137").Value = -60482

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2354.2273
$ws.Range("I40").Value = 2249.9443
$ws.Range("J40").Value = 2823.5
$ws.Range("K40").Value = 2249.9443
$ws.Range("L40").Value = 2823.5
$ws.Range("M40").Value = -2113.9443
$ws.Range("N40").Value = -3095.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 6667836
